$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.725.01"
$ws.Range("E2").Value = "  +1.87%  "

$ws.Range("D3").Value = "1.866.55"
$ws.Range("E3").Value = "  +0.57%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.52"
$ws.Range("E5").Value = "  +1.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4711"
$ws.Range("E7").Value = "  -0.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2753"
$ws.Range("E8").Value = "  +0.72%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06368"
$ws.Range("E9").Value = "  -0.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.71"
$ws.Range("E10").Value = "  +9.12%  "

$ws.Range("D11").Value = "1.870.50"
$ws.Range("E11").Value = "  +0.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07456"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.965"
$ws.Range("E13").Value = "  -1.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "85.06"
$ws.Range("E14").Value = "  -0.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6340"
$ws.Range("E15").Value = "  +0.72%  "

$ws.Range("D16").Value = "30.710.46"
$ws.Range("E16").Value = "  +2.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "243.38"
$ws.Range("E17").Value = "  +5.70%  "

$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.84"
$ws.Range("E19").Value = "  +1.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007389"
$ws.Range("E20").Value = "  +1.00%  "

$ws.Range("E21").Value = "  +0.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.998"
$ws.Range("E22").Value = "  -0.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.063"
$ws.Range("E23").Value = "  +1.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.379"
$ws.Range("E24").Value = "  +1.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.23"
$ws.Range("E25").Value = "  -0.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.21"
$ws.Range("E26").Value = "  +2.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.898"
$ws.Range("E27").Value = "  +1.06%  "

$ws.Range("E28").Value = "  +0.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.383"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.079"
$ws.Range("E30").Value = "  -1.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.866"
$ws.Range("E31").Value = "  -1.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04958"
$ws.Range("E32").Value = "  +1.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.154"
$ws.Range("E33").Value = "  +1.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7080"
$ws.Range("E34").Value = "  -1.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.713"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01912"
$ws.Range("E36").Value = "  +1.94%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.690"
$ws.Range("E37").Value = "  +1.92%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8850"
$ws.Range("E38").Value = "  -1.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.005"
$ws.Range("E39").Value = "  +2.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "105.33"
$ws.Range("E40").Value = "  -0.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("E41").Value = "  +0.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.587"
$ws.Range("E42").Value = "  +1.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4101"
$ws.Range("E43").Value = "  +0.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.76"
$ws.Range("E44").Value = "  +7.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.284"
$ws.Range("E45").Value = "  +3.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1219"
$ws.Range("E46").Value = "  +2.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.733"
$ws.Range("E47").Value = "  +0.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "33.76"
$ws.Range("E48").Value = "  +1.82%  "

$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.379"
$ws.Range("E50").Value = "  -1.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3701"
$ws.Range("E51").Value = "  +0.26%  "
